$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the edited cells keep their Text format so values are stored as strings,
# matching the inline-string cell type used in the original workbook (D=Price, E=Volume(1h)).
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "291.49"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "-3.05%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "30.67"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "-6.21%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "4.947"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "-0.16%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.07213"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "-5.97%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.805"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "-8.25%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "7.688"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "-1.82%"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.761"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "-1.02%"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.8987"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "-2.30%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1652"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "-5.64%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07710"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "-0.93%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07973"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "-6.85%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.03030"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-4.91%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.1001"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "0.05%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.001499"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "-0.91%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.005665"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "-3.81%"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.470"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "0.28%"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "2.084"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "-3.19%"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "-1.01%"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "-2.05%"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.035"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "-5.99%"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "19.74%"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.04494"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "-0.65%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.001216"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "-0.53%"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.004017"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "-8.78%"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0001251"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "-0.09%"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01586"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "-6.58%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.04403"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "-6.21%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007294"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "-2.75%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.01008"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1306"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-3.28%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.002015"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-13.67%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.009514"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "-9.06%"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00005987"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "-4.29%"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "-0.09%"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.246"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "172.74%"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.002999"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "-3.40%"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "-0.09%"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "-0.09%"
